# Add 2022-Q4 data:
#  - insert a new worksheet "2022-Q4" right after "总计" (and before "2022-Q3")
#  - fill it with the per-fund holdings detail for 2022-Q4
#  - update the "总计" (summary) sheet with a new leading row for 2022-Q4

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet before the current second sheet
#    ("2022-Q3"), matching the workbook tab order in the target.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Item(2)
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row (B1:H1) - bold + bordered style, same as every other quarter sheet.
$summary.Range("B1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)

$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

# Column A (row index) style - bold + bordered, same as the other quarter sheets.
$summary.Range("A3").Copy()
$q4.Range("A2:A15").PasteSpecial(-4122)

$rows = @(
    @(0,  "017787", "万家宏观择时多策略灵活配置混合C", "18.04", "93.54", "7.12", "1.2844", 7),
    @(1,  "519191", "万家新利灵活配置混合",             "9.85",  "92.71", "6.94", "0.6836", 5),
    @(2,  "161724", "招商中证煤炭等权指数（LOF）A",     "17.24", "93.84", "3.10", "0.5344", 5),
    @(3,  "519185", "万家精选混合A",                   "7.89",  "93.56", "6.21", "0.4900", 10),
    @(4,  "770001", "德邦优化灵活配置混合",             "2.33",  "86.37", "4.57", "0.1065", 3),
    @(5,  "015566", "万家精选混合C",                   "1.44",  "93.56", "6.21", "0.0894", 10),
    @(6,  "005944", "工银聚福混合C",                   "3.62",  "29.83", "1.37", "0.0496", 9),
    @(7,  "013596", "招商中证煤炭等权指数（LOF）C",     "1.56",  "93.84", "3.10", "0.0484", 5),
    @(8,  "003132", "德邦新回报灵活配置混合",           "0.73",  "91.15", "2.73", "0.0199", 9),
    @(9,  "012977", "瑞达鑫红量化6个月持有混合A",       "0.35",  "94.66", "4.99", "0.0175", 1),
    @(10, "016347", "招商中证煤炭等权指数（LOF）E",     "0.20",  "93.84", "3.10", "0.0062", 5),
    @(11, "012978", "瑞达鑫红量化6个月持有混合C",       "0.09",  "94.66", "4.99", "0.0045", 1),
    @(12, "005943", "工银聚福混合A",                   "0.06",  "29.83", "1.37", "0.0008", 9),
    @(13, "519212", "万家宏观择时多策略灵活配置混合A",   "0.00",  "93.54", "7.12", $null,    7)
)

$r = 2
foreach ($row in $rows) {
    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = "'" + $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4.Cells.Item($r, 6).Value = "'" + $row[5]
    if ($row[6] -eq $null) {
        $q4.Cells.Item($r, 7).Value = 0
    } else {
        $q4.Cells.Item($r, 7).Value = "'" + $row[6]
    }
    $q4.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert a new row 2 with the 2022-Q4
#    totals, shifting the rest of the table down by one row.
# ---------------------------------------------------------------------------
$summary.Rows.Item(2).Insert()

# The inserted row borrows A3's style (bold+border) for column A and clears
# the inherited header formatting on B:D so it matches every other data row.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B2:D2").ClearFormats()

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 14
$summary.Cells.Item(2, 4).Value = 3.34
